$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete the attendance record: fill in the missing "is present" value for
# row 11 (D11), and append rows 12-14 with the new attendance entries.
#
# The sheet stores everything (including numbers and dates) as plain text in
# "General" format (style index 1). When Excel is handed a numeric- or
# date-looking string it auto-converts it to a real number/date serial and
# mints a brand new style, so for every new cell we briefly force Text
# format before assigning the value, then switch the format back to
# General - the cell keeps holding the literal text we typed but regains
# the same style/format as its neighbours (style index 1).

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
}

# Row 11 - add the missing "is present" entry
Set-TextValue $ws.Range("D11") "yes"

# Row 12 - add Day Count, Date and is present
Set-TextValue $ws.Range("B12") "10"
Set-TextValue $ws.Range("C12") "07-02-2026"
Set-TextValue $ws.Range("D12") "yes"

# Row 13 - new row
Set-TextValue $ws.Range("A13") "12"
Set-TextValue $ws.Range("B13") "null"
Set-TextValue $ws.Range("C13") "08-02-2026"
Set-TextValue $ws.Range("D13") "sunday"

# Row 14 - new row (no "is present" value yet)
Set-TextValue $ws.Range("A14") "13"
Set-TextValue $ws.Range("B14") "11"
Set-TextValue $ws.Range("C14") "09-02-2026"
